$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 610; every subsequent
# record (old rows 610-659) shifts down by one row (to 611-660).
$ws.Rows.Item(610).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A610").Value = 3
$ws.Range("B610").Value = "Femacal de La Calera"
$ws.Range("C610").Value = "Coquimbo"
$ws.Range("D610").Value = 45223
$ws.Range("E610").Value = 5
$ws.Range("F610").Value = 100112040
$ws.Range("G610").Value = "Cilantro"
$ws.Range("H610").Value = "Sin especificar"
$ws.Range("I610").Value = "Primera"
$ws.Range("J610").Value = 220
$ws.Range("K610").Value = 3500
$ws.Range("L610").Value = 4000
$ws.Range("M610").Value = 3750
$ws.Range("N610").Value = "$/docena de atados (3 kilos)"
$ws.Range("O610").Value = "Provincia de Quillota"
$ws.Range("P610").Value = 1250
$ws.Range("Q610").Value = 3
$ws.Range("R610").Value = "Hortaliza"
